$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add 5 new rows (32-36) below the existing data (which ends at row 31) ---
# We copy whole rows from existing, similarly-formatted rows so the new rows
# inherit the same (mostly invisible) cell styles used throughout the table,
# then overwrite the actual values/formulas to match the new data.

# Row 32: template = row 22 (style s=15/12 pattern, N has s="12")
$ws.Rows(22).Copy()
$ws.Rows(32).Insert(-4121)

# Row 33: template = row 29 (style s=15/12 pattern, N has no explicit style)
$ws.Rows(29).Copy()
$ws.Rows(33).Insert(-4121)

# Row 34: template = row 30 (style s=15/12 pattern, N has formula, no explicit style)
$ws.Rows(30).Copy()
$ws.Rows(34).Insert(-4121)

# Row 35: template = row 22 again (N has s="12")
$ws.Rows(22).Copy()
$ws.Rows(35).Insert(-4121)

# Row 36: template = row 29 again (N has no explicit style); I-L will be cleared
$ws.Rows(29).Copy()
$ws.Rows(36).Insert(-4121)

# --- Fill in the actual values for the new rows ---

# Row 32
$ws.Range("D32").Value = 0
$ws.Range("E32").Value = "UP"
$ws.Range("F32").Value = "CAP_BND"
$ws.Range("I32").Value = 4
$ws.Range("J32").Value = 4
$ws.Range("K32").Value = 4
$ws.Range("L32").Value = 4
$ws.Range("M32").Value = 4
$ws.Range("N32").Value = "TB_SUPCO2_DKISLBH_DKE_01"

# Row 33
$ws.Range("D33").Value = 0
$ws.Range("E33").Value = "UP"
$ws.Range("F33").Value = "CAP_BND"
$ws.Range("I33").Value = 4
$ws.Range("J33").Value = 4
$ws.Range("K33").Value = 4
$ws.Range("L33").Value = 4
$ws.Range("M33").Value = 4
$ws.Range("N33").Value = "TB_H2_MAR_DKW_01"

# Row 34
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = "UP"
$ws.Range("F34").Value = "CAP_BND"
$ws.Range("I34").Value = 4
$ws.Range("J34").Value = 4
$ws.Range("K34").Value = 4
$ws.Range("L34").Value = 4
$ws.Range("M34").Value = 4
$ws.Range("N34").Formula = "=N33"

# Row 35
$ws.Range("D35").Value = 0
$ws.Range("E35").Value = "UP"
$ws.Range("F35").Value = "CAP_BND"
$ws.Range("I35").Value = 4
$ws.Range("J35").Value = 4
$ws.Range("K35").Value = 4
$ws.Range("L35").Value = 4
$ws.Range("M35").Value = 4
$ws.Range("N35").Value = "TB_SUPCO2_DKISLBH_DKE_02"

# Row 36: lower bound on PV in MAR
$ws.Range("I36:L36").Clear()
$ws.Range("D36").Value = 2030
$ws.Range("E36").Value = "LO"
$ws.Range("F36").Value = "FLO_BND"
$ws.Range("M36").Value = 1
$ws.Range("N36").Value = "TB_H2_MAR_DKW_01"

# --- Update the view so the active cell / scroll position matches ---
$ws.Range("N39").Select()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 4
